# libros.xlsx update
# "Terminados Santiago, 1 Pedro y 2 Pedro y libros integrados en archivo principal"
#
# - Mark Ezequiel (row 27) as in-progress ("*")
# - Mark Santiago, 1 Pedro and 2 Pedro (rows 60-62) as finished (1)
# - Extend the filtered/used range from row 67 to row 70 (autofilter + the
#   hidden _xlnm._FilterDatabase name) and replicate the extra
#   _FilterDatabase_* defined name that LibreOffice/Excel appends on save
# - Refresh the view: scroll back to the top and select E13

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Progress column (E) updates -------------------------------------------------

# Ezequiel (Ezekiel) is now flagged as "in progress"
$ws.Range("E27").Value2 = "*"

# Santiago, 1 Pedro, 2 Pedro are now finished -> numeric 1
$ws.Range("E60").Value2 = 1
$ws.Range("E61").Value2 = 1
$ws.Range("E62").Value2 = 1

# --- Defined names -----------------------------------------------------------

# The hidden _xlnm._FilterDatabase name now covers the extra summary rows
$wb.Names.Item(1).RefersTo = "=Sheet1!`$B`$1:`$I`$70"

# A further _FilterDatabase_* copy gets appended (matches the pattern of the
# existing ones already in the workbook)
$ws.Names.Add("_xlnm._FilterDatabase_0_0_0_0_0_0", "=Sheet1!`$B`$1:`$I`$67")

# --- AutoFilter range B1:I67 -> B1:I70 ---------------------------------------

$ws.Range("B1:I67").AutoFilter()
$ws.Range("B1:I70").AutoFilter()

# --- View state ---------------------------------------------------------------

$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("E13").Select()
